# Daily coding-problem tracker ("每日做题计划") update.
# Sheet1 ("新题") rows 31-33 are reshuffled/expanded into rows 31-41:
# two new days are inserted before the old row 31, the old rows 31-33
# are pushed down (row 33's text "63 dp" becomes a real number + category),
# four blank placeholder days are added, and a new day is appended at the end.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Make sure newly-touched date cells in column A inherit the same date
# number format ("m/d/yyyy", i.e. the style already used by A2:A30).
$ws.Range("A30").Copy()
$ws.Range("A31:A41").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 31: 4/2/2019, problem 799, done
$ws.Cells.Item(31, 1).Value = 43557
$ws.Cells.Item(31, 2).Value = 799
$ws.Cells.Item(31, 6).Value = "done"

# Row 32: 4/3/2019, problem 845, done
$ws.Cells.Item(32, 1).Value = 43558
$ws.Cells.Item(32, 2).Value = 845
$ws.Cells.Item(32, 6).Value = "done"

# Row 33: 4/4/2019, problem 552, done (this used to be row 31)
$ws.Cells.Item(33, 1).Value = 43559
$ws.Cells.Item(33, 2).Value = 552
$ws.Cells.Item(33, 6).Value = "done"

# Row 34: 4/5/2019, problem 542, category bfs, done
$ws.Cells.Item(34, 1).Value = 43560
$ws.Cells.Item(34, 2).Value = 542
$ws.Cells.Item(34, 3).Value = "bfs"
$ws.Cells.Item(34, 6).Value = "done"

# Rows 35-38: 4/6 - 4/9/2019, placeholders (date only, not done yet)
$ws.Cells.Item(35, 1).Value = 43561
$ws.Cells.Item(36, 1).Value = 43562
$ws.Cells.Item(37, 1).Value = 43563
$ws.Cells.Item(38, 1).Value = 43564

# Row 39: 4/10/2019, problem 837, done (this used to be row 32)
$ws.Cells.Item(39, 1).Value = 43565
$ws.Cells.Item(39, 2).Value = 837
$ws.Cells.Item(39, 6).Value = "done"

# Row 40: 4/11/2019, problem 63, category dp, done (this used to be row 33)
$ws.Cells.Item(40, 1).Value = 43566
$ws.Cells.Item(40, 2).Value = 63
$ws.Cells.Item(40, 3).Value = "dp"
$ws.Cells.Item(40, 6).Value = "done"

# Row 41: 4/12/2019, problem 372, category 数论 (number theory), done
$ws.Cells.Item(41, 1).Value = 43567
$ws.Cells.Item(41, 2).Value = 372
$ws.Cells.Item(41, 3).Value = "数论"
$ws.Cells.Item(41, 6).Value = "done"

# Match the author's final selection/scroll position.
$ws.Activate()
$ws.Range("F41").Select()
